$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on target cells so numeric-looking strings
# (e.g. "0.5230", "2.400") are not auto-coerced to numbers by Excel,
# matching the original inline-string (text) storage.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.995.83'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.647.58'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.85%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.87'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5230'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.12%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2616'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.40%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06352'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.72'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07715'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.658.89'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.439'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.76%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.872.20'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.60%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5504'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8237'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.87%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.73'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.993.24'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.47%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.725'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.16'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.42%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.313'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.32'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1244'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.391'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.95'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.419'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05928'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.56%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.426'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.407'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.648'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9921'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.400'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.56%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.754'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5627'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01604'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.866'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8580'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.52%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.002'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.026.40'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -7.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.07'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.795.22'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.95%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₈107'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.71'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.47%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.041'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.58%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.67%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4209'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.77%  '
